# ---------------------------------------------------------------------------
# "New model seems to be working."
#
# 1. imanes sheet: retitle the existing (left) table to make clear it is the
#    Gaussmeter measurement, and add a second (right) table holding the
#    field-intensity values produced by the new second model. Also narrow
#    the "Pequeño" average/stdev/SE down to the last 5 measurements.
# 2. add a new "distancias" sheet with the z1/l distance measurements.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "imanes"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("imanes")

# Re-word the title of the existing left-hand table (B2:D3 merged banner).
$ws.Range("B2").Value = "Intensidad del campo magnético medida con el Gaussmetro para cada imán."

# Clone all of the formatting (borders/alignment/number formats/merges) of
# the existing left table onto the new right-hand table before filling it
# with data, so both tables look identical.
$ws.Range("B2:D5").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2:H3").Merge() | Out-Null
$ws.Range("F4:F5").Merge() | Out-Null
$ws.Range("G4:H4").Merge() | Out-Null

# New right-hand table: field intensity as determined by the second model.
$ws.Range("F2").Value = "Intensidad del campo magnético determinada por el segundo modelo para cada imán."
$ws.Range("F4").Value = "Medición"
$ws.Range("G4").Value = "Int. campo (G)"
$ws.Range("G5").Value = "Grande"
$ws.Range("H5").Value = "Pequeño"

$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "`$3400\pm500`$"
$ws.Range("H6").Value = "`$970\pm70`$"

$ws.Range("F7").Value = 2
$ws.Range("G7").Value = "`$3300\pm400`$"
$ws.Range("H7").Value = "`$960\pm60`$"

$ws.Range("F8").Value = 3
$ws.Range("G8").Value = "`$3600\pm200`$"
$ws.Range("H8").Value = "`$1000\pm60`$"

$ws.Range("F9").Value = 4
$ws.Range("G9").Value = "`$3400\pm400`$"
$ws.Range("H9").Value = "`$1000\pm60`$"

$ws.Range("F10").Value = 5
$ws.Range("G10").Value = "`$3600\pm200`$"
$ws.Range("H10").Value = "`$940\pm80`$"

$ws.Columns.Item(7).AutoFit() | Out-Null

# An empty, underlined K1 cell was left behind in the source sheet.
$ws.Range("K1").Font.Underline = 2

# The "Pequeño" average/stdev/SE now only look at the last 5 measurements.
$ws.Range("D16").Formula = "=AVERAGE(D11:D15)"
$ws.Range("D17").Formula = "=STDEV(D11:D15)"

# Keep the view + selection pointed at the new table, matching the edit.
$ws.Range("K1").Select() | Out-Null

# ---------------------------------------------------------------------------
# New sheet "distancias"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dist = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$dist.Name = "distancias"

$dist.Range("A2").Value = "Medición"
$dist.Range("B2").Value = "Distancia (cm`$\pm`$0.05)"
$dist.Range("B3").Value = "`$z_1`$"
$dist.Range("C3").Value = "`$l`$"

$dist.Range("A2:A3").Merge() | Out-Null
$dist.Range("B2:C2").Merge() | Out-Null

$z1 = @(4.1, 4.1, 3.9, 3.9, 4.0, 4.0, 4.1, 4.0, 4.1, 4.0)
$l  = @(5.5, 5.4, 5.5, 5.5, 5.4, 5.4, 5.4, 5.3, 5.4, 5.5)

for ($i = 0; $i -lt 10; $i++) {
    $row = 4 + $i
    $dist.Range("A$row").Value = $i + 1
    $dist.Range("B$row").Value = $z1[$i]
    $dist.Range("C$row").Value = $l[$i]
}

$dist.Range("B4:C13").NumberFormat = "0.0"

$dist.Range("A14").Value = "Promedio (cm)"
$dist.Range("B14").Value = "`$4.02\pm0.02`$"
$dist.Range("C14").Value = "`$5.43\pm0.02`$"

$dist.Range("A17").Value = "Valor Promedio"
$dist.Range("B17").Formula = "=AVERAGE(B4:B13)"
$dist.Range("C17").Formula = "=AVERAGE(C4:C13)"
$dist.Range("B17:C17").Style = "Normal"

$dist.Range("A18").Value = "SE"
$dist.Range("B18").Formula = "=STDEV.P(B4:B13)/SQRT(10)"
$dist.Range("C18").Formula = "=STDEV.P(C4:C13)/SQRT(10)"
$dist.Range("B18:C18").Style = "Normal"

# Borders / alignment matching the rest of the workbook's small data tables.
$dist.Range("A2").Borders.Item(8).LineStyle = 1
$dist.Range("A2").Borders.Item(9).LineStyle = 1
$dist.Range("A2").HorizontalAlignment = -4108
$dist.Range("A2").VerticalAlignment = -4108
$dist.Range("A2").WrapText = $true

$dist.Range("B2:C2").Borders.Item(8).LineStyle = 1
$dist.Range("B2:C2").Borders.Item(9).LineStyle = 1
$dist.Range("B2:C2").HorizontalAlignment = -4108
$dist.Range("B2:C2").WrapText = $true

$dist.Range("A3:C3").Borders.Item(9).LineStyle = 1
$dist.Range("A3:C13").HorizontalAlignment = -4108

$dist.Range("A13:C13").Borders.Item(9).LineStyle = 1
$dist.Range("A14:C14").Borders.Item(9).LineStyle = 1

$dist.Range("K13").Font.Underline = 2

$dist.Columns.Item(1).ColumnWidth = 14.85546875
$dist.Columns.Item(2).ColumnWidth = 13.85546875
$dist.Columns.Item(3).ColumnWidth = 13.85546875

$dist.Range("B15").Select() | Out-Null

$ws.Activate()
